$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("MeV.Fib_Unknown.8").Name = "MeV.FibUnknown.8"
$wb.Worksheets.Item("MeV.Endothelial_Injury.4").Name = "MeV.EndothelialInjury.4"
$wb.Worksheets.Item("MeV.Immune_doublets.0").Name = "MeV.ImmuneDoublets.0"
$wb.Worksheets.Item("MeV.Low_Quality.0").Name = "MeV.LowQuality.0"
$wb.Worksheets.Item("MeV.Proliferative_Fibr.0").Name = "MeV.FibProlif.0"
$wb.Worksheets.Item("MeV.Epithelial_ECad.0").Name = "MeV.EpithelialECad.0"
$wb.Worksheets.Item("MeV.Fib_CD34.7").Name = "MeV.FibCD34.7"
